$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells in column D hold text-formatted numbers (e.g. "1.00", "11.60")
# that must keep trailing zeros / grouping dots; force text format before
# writing so Excel does not silently coerce them to numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.914.55'
$ws.Range("E2").Value = '  -2.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.170.20'
$ws.Range("E3").Value = '  -2.78%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.57'
$ws.Range("E5").Value = '  -2.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.616'
$ws.Range("E6").Value = '  -1.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.96'
$ws.Range("E7").Value = '  -7.93%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.564'
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.96'
$ws.Range("E10").Value = '  -0.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0922'
$ws.Range("E11").Value = '  -5.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.42'
$ws.Range("E12").Value = '  -16.02%  '
$ws.Range("E13").Value = '  -1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.89'
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.489.27'
$ws.Range("E15").Value = '  -2.71%  '
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.18'
$ws.Range("E17").Value = '  -5.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.157.10'
$ws.Range("E18").Value = '  -3.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '40.853.27'
$ws.Range("E19").Value = '  -2.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0938'
$ws.Range("E20").Value = '  -3.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.08'
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.33'
$ws.Range("E22").Value = '  -2.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '228.88'
$ws.Range("E23").Value = '  -1.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("E24").Value = '  -8.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.67'
$ws.Range("E25").Value = '  +14.09%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.72'
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.39'
$ws.Range("E28").Value = '  -4.58%  '
$ws.Range("E29").Value = '  -3.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.22'
$ws.Range("E30").Value = '  -1.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.15'
$ws.Range("E31").Value = '  -2.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.118'
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0735'
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("E35").Value = '  -3.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.55'
$ws.Range("E36").Value = '  -2.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.35'
$ws.Range("E37").Value = '  -4.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.01'
$ws.Range("E38").Value = '  -2.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0298'
$ws.Range("E39").Value = '  +4.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.17'
$ws.Range("E40").Value = '  -5.24%  '
$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.45'
$ws.Range("E41").Value = '  -9.90%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.60'
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.81'
$ws.Range("E43").Value = '  -14.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.79'
$ws.Range("E44").Value = '  -6.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.189'
$ws.Range("E45").Value = '  -11.37%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.45'
$ws.Range("E46").Value = '  -4.69%  '
$ws.Range("B47").Value = 'BinanceUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0990'
$ws.Range("E48").Value = '  -2.93%  '
$ws.Range("E49").Value = '  -0.83%  '
$ws.Range("E50").Value = '  -3.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.69'
$ws.Range("E51").Value = '  -1.26%  '

Write-Output "Applied 98 cell updates"
